$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '43.826.96'
$ws.Range('E2').Value = '  -0.30%  '
$ws.Range('D3').Value = '2.317.89'
$ws.Range('E3').Value = '  +3.25%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '97.16'
$ws.Range('E5').Value = '  +5.17%  '
$ws.Range('D6').Value = '272.42'
$ws.Range('E6').Value = '  +0.37%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '0.626'
$ws.Range('E9').Value = '  -0.38%  '
$ws.Range('D10').Value = '45.37'
$ws.Range('E10').Value = '  -0.89%  '
$ws.Range('D11').Value = '0.0953'
$ws.Range('E11').Value = '  -1.73%  '
$ws.Range('D12').Value = '8.04'
$ws.Range('E12').Value = '  -3.37%  '
$ws.Range('D13').Value = '0.105'
$ws.Range('E13').Value = '  -0.11%  '
$ws.Range('D14').Value = '2.657.61'
$ws.Range('E14').Value = '  +2.93%  '
$ws.Range('D15').Value = '15.54'
$ws.Range('E15').Value = '  +2.52%  '
$ws.Range('D16').Value = '0.875'
$ws.Range('E16').Value = '  +7.96%  '
$ws.Range('D17').Value = '2.314.91'
$ws.Range('E17').Value = '  +2.98%  '
$ws.Range('D18').Value = '43.771.15'
$ws.Range('E18').Value = '  -0.30%  '
$ws.Range('E19').Value = '  +4.41%  '
$ws.Range('D20').Value = '6.42'
$ws.Range('E20').Value = '  +4.82%  '
$ws.Range('D21').Value = '73.41'
$ws.Range('E21').Value = '  +3.60%  '
$ws.Range('D22').Value = '239.80'
$ws.Range('E22').Value = '  +2.10%  '
$ws.Range('E23').Value = '  -3.00%  '
$ws.Range('D24').Value = '9.44'
$ws.Range('E25').Value = '  -0.08%  '
$ws.Range('D26').Value = '2.54'
$ws.Range('E26').Value = '  +1.24%  '
$ws.Range('D27').Value = '11.39'
$ws.Range('E27').Value = '  -0.71%  '
$ws.Range('E28').Value = '  -1.28%  '
$ws.Range('E29').Value = '  +1.39%  '
$ws.Range('D30').Value = '38.33'
$ws.Range('E30').Value = '  -7.42%  '
$ws.Range('D31').Value = '22.41'
$ws.Range('E31').Value = '  +6.69%  '
$ws.Range('D32').Value = '174.91'
$ws.Range('E32').Value = '  +1.24%  '
$ws.Range('D33').Value = '0.0918'
$ws.Range('E33').Value = '  -0.09%  '
$ws.Range('E34').Value = '  -0.11%  '
$ws.Range('E35').Value = '  +2.68%  '
$ws.Range('D36').Value = '0.0367'
$ws.Range('E36').Value = '  +4.52%  '
$ws.Range('D37').Value = '0.110'
$ws.Range('E37').Value = '  -3.75%  '
$ws.Range('E38').Value = '  +3.18%  '
$ws.Range('D39').Value = '3.38'
$ws.Range('E39').Value = '  -4.58%  '
$ws.Range('D40').Value = '0.245'
$ws.Range('E40').Value = '  +8.03%  '
$ws.Range('E41').Value = '  +10.69%  '
$ws.Range('D42').Value = '1.42'
$ws.Range('E42').Value = '  +22.99%  '
$ws.Range('D43').Value = '12.34'
$ws.Range('E43').Value = '  -5.43%  '
$ws.Range('D44').Value = '62.85'
$ws.Range('E44').Value = '  -1.66%  '
$ws.Range('D45').Value = '9.22'
$ws.Range('E45').Value = '  +10.10%  '
$ws.Range('E46').Value = '  -0.56%  '
$ws.Range('E47').Value = '  +3.39%  '
$ws.Range('D48').Value = '100.54'
$ws.Range('E48').Value = '  +0.21%  '
$ws.Range('E49').Value = '  +0.74%  '
$ws.Range('D50').Value = '0.193'
$ws.Range('E50').Value = '  +17.21%  '
$ws.Range('D51').Value = '2.544.63'
$ws.Range('E51').Value = '  +3.09%  '
